# University Program view, add, update, import, export, bulk update
# Rebuild header row with the new field set/order and drop the unused
# exam-specific columns (ielts/toefl/pte/duolingo/gre/gmat/sat, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "course_name"
$ws.Range("B1").Value = "course_category_id"
$ws.Range("C1").Value = "specialization_id"
$ws.Range("D1").Value = "level"
$ws.Range("E1").Value = "duration"
$ws.Range("F1").Value = "study_mode"
$ws.Range("G1").Value = "intake"
$ws.Range("H1").Value = "application_deadline"
$ws.Range("I1").Value = "tution_fee"
$ws.Range("J1").Value = "overview"
$ws.Range("K1").Value = "entry_requirement"
$ws.Range("L1").Value = "exam_required"
$ws.Range("M1").Value = "mode_of_instruction"
$ws.Range("N1").Value = "scholarship_info"

# Drop the now-unused trailing columns (old sheet went out to S1)
$ws.Range("O1:S1").ClearContents()

# Widen the application_deadline column like the source workbook
$ws.Columns.Item(8).ColumnWidth = 26.5

# Match the author's final active cell/selection
[void]$ws.Range("O1").Select()
